# Restore the "Cash" calculation values that were accidentally dropped.
# This updates the Weekly Performance (%), Performance (%) and Value (£)
# columns for every holding row (2-18) plus the Performance (%) and
# Value (£) totals on the "Cash" row (19).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Weekly_Update")

$ws.Range("B2").Value = -1.98
$ws.Range("C2").Value = -37.59
$ws.Range("D2").Value = 648.2265110778809

$ws.Range("B3").Value = -0.21
$ws.Range("C3").Value = 10.83
$ws.Range("D3").Value = 1234.248243933589

$ws.Range("B4").Value = -2.09
$ws.Range("C4").Value = -41.67
$ws.Range("D4").Value = 238.0061077875848

$ws.Range("B5").Value = -2.47
$ws.Range("C5").Value = -21.36
$ws.Range("D5").Value = 201.1710690622063

$ws.Range("B6").Value = 2
$ws.Range("C6").Value = 132.36
$ws.Range("D6").Value = 1304.761030125109

$ws.Range("B7").Value = 3.5
$ws.Range("C7").Value = -29.03
$ws.Range("D7").Value = 514.0410280970626

$ws.Range("B8").Value = -0.1
$ws.Range("C8").Value = -35.37
$ws.Range("D8").Value = 297.2563085738727

$ws.Range("B9").Value = 6.5
$ws.Range("C9").Value = -67.02
$ws.Range("D9").Value = 398.0999908447266

$ws.Range("B10").Value = -4.08
$ws.Range("C10").Value = -5.51
$ws.Range("D10").Value = 571.5

$ws.Range("B11").Value = 0.08
$ws.Range("C11").Value = 2.04
$ws.Range("D11").Value = 521.6324033409119

$ws.Range("B12").Value = 0.22
$ws.Range("C12").Value = 47.72
$ws.Range("D12").Value = 882.6889677944183

$ws.Range("B13").Value = 7.79
$ws.Range("C13").Value = -6.75
$ws.Range("D13").Value = 1384.415729758644

$ws.Range("B14").Value = -2.16
$ws.Range("C14").Value = 28.04
$ws.Range("D14").Value = 895.50735367495

$ws.Range("B15").Value = 0.64
$ws.Range("C15").Value = -29.46
$ws.Range("D15").Value = 489.5052869503047

$ws.Range("B16").Value = 0.7
$ws.Range("C16").Value = 24.42
$ws.Range("D16").Value = 557.9988966590881

$ws.Range("B17").Value = -2.08
$ws.Range("C17").Value = -74.22
$ws.Range("D17").Value = 98.28432993479093

$ws.Range("B18").Value = 1.11
$ws.Range("C18").Value = 2.37
$ws.Range("D18").Value = 10237.34325761514

$ws.Range("C19").Value = 1323.367396896345
$ws.Range("D19").Value = 912.6268620576598
